$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53, shifting existing rows 53-172 down to 54-173.
$ws.Rows("53:53").Insert()

# Populate the new row 53 with its data: the columns that stay constant
# throughout the table (A,B,C,E,F,G,H,I,O,R) are copied from the row that
# used to occupy position 53 (now at 54), and the columns that actually
# differ for this new record (D,J,K,L,M,N,P,Q) get their new values.
$ws.Range("A53").Value = 11
$ws.Range("B53").Value = "Vega Monumental Concepción"
$ws.Range("C53").Value = "Bíobío"
$ws.Range("D53").Value = 45125
$ws.Range("E53").Value = 8
$ws.Range("F53").Value = 100112001
$ws.Range("G53").Value = "Berenjena"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 170
$ws.Range("K53").Value = 7000
$ws.Range("L53").Value = 8000
$ws.Range("M53").Value = 7471
$ws.Range("N53").Value = "$/caja 50 unidades"
$ws.Range("O53").Value = "Región de Arica y Parinacota"
$ws.Range("P53").Value = 149
$ws.Range("Q53").Value = 50
$ws.Range("R53").Value = "Hortaliza"
